$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add context to final report: update revenue figures for Australia, New
# Zealand, Canada, USA and Japan rows.
$ws.Range("C3").Value = 22
$ws.Range("C4").Value = 67
$ws.Range("C7").Value = 71
$ws.Range("C8").Value = 66
$ws.Range("C9").Value = 3

# Leave the cursor/selection on the last edited cell.
$ws.Range("C9").Select()

# Match the updated workbook window placement/size from the saved file.
$excel.ActiveWindow.Left = 19200
$excel.ActiveWindow.Top = 0
$excel.ActiveWindow.Width = 19200
$excel.ActiveWindow.Height = 21000

$wb.Save()
